$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.944.00'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '1.639.70'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  +0.86%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.00'
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("E7").Value = '  +0.82%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.255'
$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0639'
$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("E10").Value = '  -0.45%  '

$ws.Range("E11").Value = '  +0.88%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.678.65'
$ws.Range("E12").Value = '  +3.05%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.867.01'
$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.25'
$ws.Range("E14").Value = '  +0.15%  '

$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("E16").Value = '  +0.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.66'
$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("D18").Value = '25.964.32'
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("E19").Value = '  +0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.37'
$ws.Range("E20").Value = '  +0.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.37'
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.94'
$ws.Range("E22").Value = '  -0.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.29'
$ws.Range("E23").Value = '  -0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.22'
$ws.Range("E24").Value = '  +1.60%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  +1.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.128'
$ws.Range("E27").Value = '  +2.16%  '

$ws.Range("E28").Value = '  -0.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.49'
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("E31").Value = '  +0.73%  '

$ws.Range("E32").Value = '  -1.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("E34").Value = '  -2.60%  '

$ws.Range("E35").Value = '  +1.73%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.904'
$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("D37").Value = '1.140.09'
$ws.Range("E37").Value = '  +0.29%  '

$ws.Range("E38").Value = '  -0.56%  '

$ws.Range("E39").Value = '  -1.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0157'
$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.36'
$ws.Range("E41").Value = '  -0.73%  '

$ws.Range("E42").Value = '  +1.49%  '

$ws.Range("E43").Value = '  -2.70%  '

$ws.Range("D44").Value = '1.776.51'
$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("D45").Value = '0.0₆0116'
$ws.Range("E45").Value = '  +8.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.58'
$ws.Range("E46").Value = '  +1.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0529'
$ws.Range("E47").Value = '  +2.54%  '

$ws.Range("E48").Value = '  +0.10%  '

$ws.Range("E49").Value = '  -0.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.64'
$ws.Range("E50").Value = '  -0.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0965'
$ws.Range("E51").Value = '  -0.75%  '
